$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test-Cases")

$ws.Range("J5").Value = "sdtdft"
$ws.Range("J7").Value = "fyfy"

$ws.Range("J7").Select()
